# Updated symbol list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row on the sheet. Values are stored as plain text in the workbook
# (matching the existing inlineStr cells), so each literal is written with a
# leading apostrophe to tell Excel to keep it as text instead of
# reinterpreting the numeric-looking string as a Number/Percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.54"
$ws.Range("E2").Value = "'0.73%"
$ws.Range("D3").Value = "'29.16"
$ws.Range("E3").Value = "'1.40%"
$ws.Range("D4").Value = "'5.287"
$ws.Range("E4").Value = "'4.80%"
$ws.Range("D5").Value = "'0.07083"
$ws.Range("E5").Value = "'5.68%"
$ws.Range("D6").Value = "'7.455"
$ws.Range("E6").Value = "'1.68%"
$ws.Range("D7").Value = "'3.561"
$ws.Range("E7").Value = "'5.15%"
$ws.Range("D8").Value = "'1.392"
$ws.Range("E8").Value = "'1.51%"
$ws.Range("D9").Value = "'0.9041"
$ws.Range("E9").Value = "'-3.82%"
$ws.Range("D10").Value = "'0.1603"
$ws.Range("E10").Value = "'2.77%"
$ws.Range("D11").Value = "'0.07545"
$ws.Range("E11").Value = "'11.48%"
$ws.Range("D12").Value = "'0.07713"
$ws.Range("E12").Value = "'1.88%"
$ws.Range("D13").Value = "'0.02912"
$ws.Range("E13").Value = "'-1.45%"
$ws.Range("E14").Value = "'0.38%"
$ws.Range("D15").Value = "'0.001603"
$ws.Range("E15").Value = "'0.83%"
$ws.Range("D16").Value = "'0.0006510"
$ws.Range("E16").Value = "'0.74%"
$ws.Range("D17").Value = "'0.006371"
$ws.Range("E17").Value = "'-3.01%"
$ws.Range("D18").Value = "'3.494"
$ws.Range("E18").Value = "'1.35%"
$ws.Range("E19").Value = "'-0.62%"
$ws.Range("E20").Value = "'-0.01%"
$ws.Range("D22").Value = "'4.008"
$ws.Range("E22").Value = "'-1.32%"
$ws.Range("E23").Value = "'3.01%"
$ws.Range("D24").Value = "'0.04524"
$ws.Range("E24").Value = "'0.81%"
$ws.Range("D25").Value = "'0.001210"
$ws.Range("E25").Value = "'2.56%"
$ws.Range("D26").Value = "'0.004167"
$ws.Range("E26").Value = "'-7.14%"
$ws.Range("E27").Value = "'-6.16%"
$ws.Range("E28").Value = "'3.12%"
$ws.Range("D40").Value = "'0.04378"
$ws.Range("E40").Value = "'4.12%"
$ws.Range("D41").Value = "'0.006987"
$ws.Range("E41").Value = "'4.15%"
$ws.Range("E42").Value = "'-0.23%"
$ws.Range("E43").Value = "'2.74%"
$ws.Range("D44").Value = "'0.01171"
$ws.Range("E44").Value = "'-4.83%"
$ws.Range("D45").Value = "'0.00005859"
$ws.Range("E45").Value = "'3.73%"
$ws.Range("E47").Value = "'-0.57%"
